$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a 1-row array of 36 values (F..AO) into the given row number
function Set-RowValues {
    param($RowNum, $Values)
    $arr = New-Object "object[,]" 1,36
    for ($i = 0; $i -lt 36; $i++) { $arr[0,$i] = $Values[$i] }
    $rangeStr = "F$RowNum`:AO$RowNum"
    $ws.Range($rangeStr).Value = $arr
}

$row2 = @(1.62,1.64,5.6,5.7,4.7,4.8,1.31,1.04,5.5,1.2,2.56,1.62,1.61,2.56,1.69,2.36,1.21,2.56,27,26,50,140,11.5,11,21,65,11.5,9.6,18.5,60,16,15,27,80,7.2,60)
$row3 = @(27,55,1.07,1.09,16,22,1.09,1.01,15,1.04,6.4,1.15,3.2,1.37,2,1.81,12,1.02,100,36,20,14,1000,55,27,21,1000,1000,990,55,1000,1000,1000,1000,1000,1.75)
$row4 = @(2.18,2.42,4,4.8,2.84,3.2,1.61,1.12,2.34,1.6,1.44,2.82,1.15,6,2.24,1.66,1.27,1.7,14,20,1000,1000,7,7.6,60,1000,34,23,1000,1000,1000,95,1000,1000,1000,1000)
$row5 = @(2.24,2.34,3,3.2,3.95,4.4,1.01,1.03,6.2,1.15,2.8,1.46,1.75,2.12,1.5,2.86,1.45,1.74,80,40,85,65,34,11.5,15,1000,22,13,14.5,80,120,40,55,580,10,15.5)
$row6 = @(2.54,2.7,3.45,3.9,2.72,3.05,1.66,1.15,2.2,1.66,1.39,3,1.14,6.4,2.28,1.63,1.36,1.58,7.6,10.5,28,110,7.2,7.2,22,190,16.5,15.5,34,130,46,48,95,1000,70,140)
$row7 = @(1.45,1.5,7.6,9,4.6,5.1,1.28,1.05,4.1,1.26,2.08,1.76,1.44,2.9,1.9,1.89,1.12,2.96,20,65,1000,300,8.8,11,80,150,9.4,10.5,46,130,14,17,85,160,7.4,1000)
$row8 = @(2.72,2.92,2.68,2.9,3.4,3.6,1.01,1.07,3.35,1.35,1.83,2.02,1.31,3.6,1.77,2.06,1.53,1.52,16,11,24,280,13.5,9.4,15.5,40,24,15.5,22,60,55,40,55,580,32,36)
$row9 = @(1.6,1.64,9.4,12,3.35,3.65,1.57,1.15,2.3,1.61,1.43,2.84,1.14,6.2,2.84,1.44,1.1,2.56,15,60,1000,1000,8.8,17.5,1000,1000,40,36,1000,1000,180,140,1000,1000,1000,1000)
$row10 = @(1.7,1.71,7,7.2,3.7,3.75,1.55,1.11,2.96,1.49,1.65,2.48,1.24,5,2.36,1.7,1.16,2.42,9.8,17,55,230,6.2,8.4,28,140,8,11,29,150,16,22,55,250,15,230)
$row11 = @(1.1,1.11,46,50,12.5,13,1.22,1.02,7.6,1.13,3.05,1.44,1.81,2.12,2.84,1.51,1.02,10,46,990,620,1000,12.5,990,990,1000,7.4,17.5,990,1000,7.2,16,85,590,2.94,1000)

Set-RowValues 2 $row2
Set-RowValues 3 $row3
Set-RowValues 4 $row4
Set-RowValues 5 $row5
Set-RowValues 6 $row6
Set-RowValues 7 $row7
Set-RowValues 8 $row8
Set-RowValues 9 $row9
Set-RowValues 10 $row10
Set-RowValues 11 $row11

# New row 12: League/Home/Away are plain text; Date/Time must stay text too
$ws.Range("A12").Value = "Friendly Matches"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2025-12-29"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "21:00:00"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").Value = "Tlaxcala F.C"
$ws.Range("E12").Value = "Pachuca"

$row12 = @(1.04,1000,1.04,1000,1.02,950,1.01,1.01,1.24,1.01,1.24,1.01,1.18,1.01,1.01,1.01,1.01,1.01,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)
Set-RowValues 12 $row12
